$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay text (match source "inlineStr" typing)
$textCells = @("D5", "D10", "D15", "D17", "D20", "D22", "D23", "D25", "D28", "D29", "D32", "D33", "D34", "D35", "D41", "D42", "D45", "D47", "D48")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated price / volume figures
$ws.Range('D2').Value = '26.281.26'
$ws.Range('E2').Value = '  +1.12%  '
$ws.Range('D3').Value = '1.622.65'
$ws.Range('E3').Value = '  +1.39%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '212.57'
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('E7').Value = '  +0.35%  '
$ws.Range('E8').Value = '  +0.78%  '
$ws.Range('E9').Value = '  +0.26%  '
$ws.Range('D10').Value = '18.92'
$ws.Range('E10').Value = '  +4.44%  '
$ws.Range('E11').Value = '  +0.63%  '
$ws.Range('D12').Value = '1.849.57'
$ws.Range('E12').Value = '  +1.43%  '
$ws.Range('D13').Value = '1.638.84'
$ws.Range('E13').Value = '  +2.31%  '
$ws.Range('E14').Value = '  +0.52%  '
$ws.Range('D15').Value = '0.519'
$ws.Range('E15').Value = '  +1.01%  '
$ws.Range('D16').Value = '26.287.95'
$ws.Range('E16').Value = '  +1.13%  '
$ws.Range('D17').Value = '62.43'
$ws.Range('E17').Value = '  +3.74%  '
$ws.Range('E18').Value = '  +0.92%  '
$ws.Range('E19').Value = '  +0.10%  '
$ws.Range('D20').Value = '203.22'
$ws.Range('E20').Value = '  +0.73%  '
$ws.Range('E21').Value = '  +1.44%  '
$ws.Range('D22').Value = '9.36'
$ws.Range('E22').Value = '  +0.93%  '
$ws.Range('D23').Value = '6.03'
$ws.Range('E23').Value = '  +0.39%  '
$ws.Range('E24').Value = '  +7.72%  '
$ws.Range('D25').Value = '142.66'
$ws.Range('E25').Value = '  +0.67%  '
$ws.Range('E26').Value = '  +0.12%  '
$ws.Range('E27').Value = '  -0.15%  '
$ws.Range('D28').Value = '15.23'
$ws.Range('E28').Value = '  +0.82%  '
$ws.Range('D29').Value = '6.55'
$ws.Range('E29').Value = '  +1.91%  '
$ws.Range('E31').Value = '  +0.56%  '
$ws.Range('D32').Value = '3.18'
$ws.Range('E32').Value = '  +2.52%  '
$ws.Range('D33').Value = '2.95'
$ws.Range('E33').Value = '  -0.13%  '
$ws.Range('D34').Value = '1.50'
$ws.Range('E34').Value = '  +1.68%  '
$ws.Range('D35').Value = '2.41'
$ws.Range('E35').Value = '  +2.39%  '
$ws.Range('D36').Value = '1.173.72'
$ws.Range('E36').Value = '  +4.29%  '
$ws.Range('E37').Value = '  +1.41%  '
$ws.Range('E38').Value = '  +2.43%  '
$ws.Range('E40').Value = '  +0.01%  '
$ws.Range('D41').Value = '0.497'
$ws.Range('E41').Value = '  +1.18%  '
$ws.Range('D42').Value = '0.794'
$ws.Range('E42').Value = '  +0.99%  '
$ws.Range('E43').Value = '  +3.33%  '
$ws.Range('D44').Value = '1.761.31'
$ws.Range('E44').Value = '  +1.48%  '
$ws.Range('D45').Value = '93.47'
$ws.Range('E45').Value = '  +0.51%  '
$ws.Range('E46').Value = '  +14.54%  '
$ws.Range('D47').Value = '1.51'
$ws.Range('E47').Value = '  +0.71%  '
$ws.Range('D48').Value = '54.11'
$ws.Range('E48').Value = '  +1.09%  '
$ws.Range('E49').Value = '  +0.93%  '
$ws.Range('E50').Value = '  +0.19%  '
$ws.Range('E51').Value = '  -0.13%  '
